# Update data through 09.09.2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the tiny float corrections on the last two existing rows (1324, 1325)
$ws.Range("B1324").Value = 0.9188924280015596
$ws.Range("B1325").Value = 0.9196296105298016

$dates = @(
    "2023-08-18","2023-08-19","2023-08-20","2023-08-21","2023-08-22","2023-08-23","2023-08-24",
    "2023-08-25","2023-08-26","2023-08-27","2023-08-28","2023-08-29","2023-08-30","2023-08-31",
    "2023-09-01","2023-09-02","2023-09-03","2023-09-04","2023-09-05","2023-09-06","2023-09-07",
    "2023-09-08","2023-09-09"
)

$prices = @(
    0.9186590956221447,0.9208352491360975,0.9267152860013571,0.9173868726084622,0.9216308215777895,
    0.916416413647274,0.9249150713156021,0.9253015746700953,0.9260709106473283,0.9291095803036642,
    0.9245288456351862,0.9182388455392471,0.9140910262678169,0.9183313397383227,0.9275432091775185,
    0.922686074170466,0.9284101405876203,0.9265092605722369,0.9320291692234548,0.9321445941858783,
    0.9334941303910537,0.9327919389566486,0.9327130559698409
)

$startRow = 1326
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $dateCell = $ws.Cells.Item($row, 1)
    # Force the date-like string to be stored as plain text (shared string)
    # instead of being auto-converted into a date serial number.
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $dates[$i]
    $dateCell.ClearFormats()

    $ws.Cells.Item($row, 2).Value = $prices[$i]
}
